$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 to the new email address
$ws.Range("A2").Value = "patti.thompson@enron.com"

# A5 currently holds johnny.palmer, which is being removed entirely, so
# replace it with the value that should end up there (bobette.riner, which
# previously lived in A7).
$ws.Range("A5").Value = "bobette.riner@ipgdirect.com"

# Rows 6 and 7 (the old patti.thompson and bobette.riner rows) are no
# longer needed now that bobette.riner has moved up to A5.
$ws.Rows("6:7").Delete()
